# Update cryptos list values per upstream data refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.293.46"
$ws.Range("E2").Value = "  +0.64%  "
$ws.Range("D3").Value = "2.605.91"
$ws.Range("E3").Value = "  +0.51%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "543.28"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.99"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("E7").Value = "  +0.27%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("E10").Value = "  +1.75%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.334"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.39%  "
$ws.Range("E12").Value = "  +1.39%  "
$ws.Range("D13").Value = "3.065.37"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "59.227.37"
$ws.Range("E14").Value = "  +0.49%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "20.54"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -0.06%  "
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "2.648.09"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000133"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "343.14"
$ws.Range("D18").ClearFormats()
$ws.Range("E19").Value = "  +0.65%  "
$ws.Range("E20").Value = "  -1.01%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.39"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.02%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "67.50"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.93%  "
$ws.Range("B24").Value = "Polygon"
$ws.Range("C24").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.408"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.42%  "
$ws.Range("B25").Value = "Kaspa"
$ws.Range("C25").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.166"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.07%  "
$ws.Range("E26").Value = "  +0.24%  "
$ws.Range("E27").Value = "  +1.09%  "
$ws.Range("E28").Value = "  +0.10%  "
$ws.Range("E29").Value = "  +1.31%  "
$ws.Range("E30").Value = "  +8.41%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "5.79"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "18.74"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "149.90"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.65%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "3.97"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -1.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "37.12"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +2.14%  "
$ws.Range("E36").Value = "  -1.78%  "
$ws.Range("E37").Value = "  +0.34%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.832"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("E39").Value = "  -1.86%  "
$ws.Range("E40").Value = "  +0.00%  "
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.32%  "
$ws.Range("B42").Value = "Bittensor"
$ws.Range("C42").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "276.54"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -0.34%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.596"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.97%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.73"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.06%  "
$ws.Range("E45").Value = "  -0.03%  "
$ws.Range("E46").Value = "  +0.53%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0223"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.08%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Value = "1.942.85"
$ws.Range("E48").Value = "  -2.20%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "18.38"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +1.72%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.50"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "111.07"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.76%  "
